# Perbaikan Antrian Device Presensi
# Update attendance statuses for specific students and refresh the summary counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update individual status cells (column D)
$ws.Range("D2").Value = "Alpha"
$ws.Range("D3").Value = "Hadir"
$ws.Range("D4").Value = "Hadir"
$ws.Range("D7").Value = "Hadir"

# Update the "Ringkasan Kehadiran" (attendance summary) totals
$ws.Range("A37").Value = "Hadir: 32"
$ws.Range("A38").Value = "Izin: 0"
$ws.Range("A39").Value = "Sakit: 0"
